$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new row 13 data for the new minigame ("a new vote game")
# Shared-string pool order matters: new unique strings must be introduced
# in the order t10, GameButton10, 村长选举 to match the target indices.
$ws.Range("A13").Value = 17000010
$ws.Range("H13").Value = "t10"
$ws.Range("G13").Value = "GameButton10"
$ws.Range("B13").Value = "村长选举"
$ws.Range("C13").Value = 60
$ws.Range("D13").Value = 70
$ws.Range("E13").Value = 80
$ws.Range("F13").Value = 1109

# Update selection to match diff (B13 active cell)
$ws.Range("B13").Select()

# Resize the table to include the new row
$table = $ws.ListObjects.Item(1)
$table.Resize($ws.Range("A3:H13"))
